# Applies the "falas do matheus alterada" edit:
#  - Adds an intro paragraph + blank line at the very top of the document.
#  - Adds a new bold "INSTALAÇÃO" heading followed by an installation
#    instructions paragraph (carries the _GoBack bookmark now).
#  - Adds two paragraphs (one bold-formatted, empty) before "PROCESSOS DE
#    SUPORTE".
#  - Rewrites the "PROCESSOS DE SUPORTE" body paragraph, fixing a long list
#    of typos/missing accents and dropping the now-unneeded w:proofErr
#    wrappers around the corrected words.
#  - Removes the old _GoBack bookmark that used to sit in the empty bold
#    paragraph right before the "FERRAMENTA" section's chat paragraph.
#  - Fixes a handful of accent typos in the closing "FERRAMENTA" paragraph.
#
# The whole body (every paragraph + the trailing sectPr) is rebuilt in one
# shot and dropped in via Range.InsertXML so paragraph/run formatting is
# expressed explicitly instead of being inherited from whatever range we
# started the edit on.

$d = $word.ActiveDocument

$newBodyXml = @'
<w:p><w:r><w:t xml:space="preserve">Criamos um sistema que mede a temperatura e umidade do ar de incubadoras neonatal este sistema consiste em exibir os dados de temperatura e umidade em gráficos, cadastro de incubadoras e de recém-nascidos </w:t></w:r></w:p><w:p/><w:p w14:paraId="5C953A8A" w14:textId="29AB7968" w:rsidR="00137CE0" w:rsidRPr="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>INSTALAÇÃO</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Conectaremos o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arduino</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> em um computador por uma porta </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>usb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, após isso iremos dar </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>install</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para instalar as dependências do sistema que server para fazer o sistema funcionar corretamente depois daremos um </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>npm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> start para ativar o servidor e enviar os dados do sensor do </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>arduino</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> para o nosso banco de dados na nuvem logo após isso deverá ser acessado o link do site, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>logar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> e verificar se o sistema está funcionando corretamente cadastrando incubadoras e recém-nascidos.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>PROCESSOS DE SUPORTE</w:t></w:r></w:p><w:p w14:paraId="63F39CE4" w14:textId="77777777" w:rsidR="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"/><w:p w14:paraId="13D6D217" w14:textId="77777777" w:rsidR="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"><w:r><w:t>-</w:t></w:r><w:r><w:t>Temos</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>também</w:t></w:r><w:r><w:t xml:space="preserve"> um processo de atendimento </w:t></w:r><w:r><w:t>e suporte aonde teremos 3 níveis</w:t></w:r><w:r><w:t xml:space="preserve"> de atendimento. Quando o cliente. Entrar em contato por telefonema ou por chat ele </w:t></w:r><w:r><w:t>será</w:t></w:r><w:r><w:t xml:space="preserve"> direcionado para o </w:t></w:r><w:r><w:t>nível</w:t></w:r><w:r><w:t xml:space="preserve"> 1 aonde ele </w:t></w:r><w:r><w:t>será</w:t></w:r><w:r><w:t xml:space="preserve"> identificado pelo atendente, o atendente </w:t></w:r><w:r><w:t>irá</w:t></w:r><w:r><w:t xml:space="preserve"> fazer algumas perguntas simples para ver se o problema </w:t></w:r><w:r><w:t>não</w:t></w:r><w:r><w:t xml:space="preserve"> possa ser algo simples e assim solucionar caso </w:t></w:r><w:r><w:t>não</w:t></w:r><w:r><w:t xml:space="preserve"> consiga ele passara o atendimento para o </w:t></w:r><w:r><w:t>nível</w:t></w:r><w:r><w:t xml:space="preserve"> dois onde </w:t></w:r><w:r><w:t>será</w:t></w:r><w:r><w:t xml:space="preserve"> um </w:t></w:r><w:r><w:t>técnico</w:t></w:r><w:r><w:t xml:space="preserve"> mais </w:t></w:r><w:r><w:t>especializado</w:t></w:r><w:r><w:t xml:space="preserve">, ele </w:t></w:r><w:r><w:t>irá</w:t></w:r><w:r><w:t xml:space="preserve"> tentar solucionar o problema remotamente caso n consiga ira marcar uma visita presencial para q assim solucione o problema, e o </w:t></w:r><w:r><w:t>nível</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>três</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>serão</w:t></w:r><w:r><w:t xml:space="preserve"> analistas com </w:t></w:r><w:r><w:t>experiência</w:t></w:r><w:r><w:t xml:space="preserve"> para </w:t></w:r><w:r><w:t>solucionar</w:t></w:r><w:r><w:t xml:space="preserve"> problemas na arquitetura do sistema analisando </w:t></w:r><w:r><w:t>tudo</w:t></w:r><w:r><w:t xml:space="preserve"> para encontrar a causa raiz e para fazer a </w:t></w:r><w:r><w:t>alteração abrira uma GMUD</w:t></w:r><w:r><w:t xml:space="preserve"> para explicar como </w:t></w:r><w:r><w:t>será feita a mudança e as</w:t></w:r><w:r><w:t>sim solucionar o problema.</w:t></w:r></w:p><w:p/><w:p w14:paraId="5349015F" w14:textId="1F3955DB" w:rsidR="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"/><w:p w14:paraId="32C68DE0" w14:textId="539160AF" w:rsidR="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>FERRAMENTA</w:t></w:r></w:p><w:p w14:paraId="0B370998" w14:textId="77777777" w:rsidR="00137CE0" w:rsidRPr="00137CE0" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p w14:paraId="5DFC622D" w14:textId="063DA797" w:rsidR="00A36828" w:rsidRDefault="00137CE0" w:rsidP="00137CE0"><w:r><w:t xml:space="preserve">Para o atendimento por chat iremos usar a ferramenta de help </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>desk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tomticket</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> aonde </w:t></w:r><w:r><w:t>terá</w:t></w:r><w:r><w:t xml:space="preserve"> um roteiro de atendimento para nossos </w:t></w:r><w:r><w:t>funcionários</w:t></w:r><w:r><w:t xml:space="preserve"> sigam e uma base de erros conhecidos para verificar se o problema do cliente </w:t></w:r><w:r><w:t>já</w:t></w:r><w:r><w:t xml:space="preserve"> ocorreu para pode solucionar facilmente</w:t></w:r></w:p><w:sectPr w:rsidR="00A36828"><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1417" w:right="1701" w:bottom="1417" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr>
'@

$d.Content.InsertXML($newBodyXml)
